$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1: title slide date textbox - merge 3 runs ("June " / "15th, " / "2015")
# into a single run "June 15th, 2015" (keep the first run's formatting).
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$dateBox = $s1.Shapes.Item(3)
$dateRange = $dateBox.TextFrame.TextRange
# Force through an intermediate value so the engine actually rewrites the
# backing runs instead of treating an identical concatenation as a no-op.
$dateRange.Text = "TEMP_PLACEHOLDER"
$dateBox.TextFrame.TextRange.Text = "June 15th, 2015"

# ---------------------------------------------------------------------------
# Slide 3: title "Apache " + "Flink's" + " Type System" -> "Type System and Keys"
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$title3 = $s3.Shapes.Item(1)
$title3Range = $title3.TextFrame.TextRange
$title3Range.Text = "TEMP_PLACEHOLDER"
$title3.TextFrame.TextRange.Text = "Type System and Keys"

# ---------------------------------------------------------------------------
# Slide 4: content placeholder
#   - paragraph 1: "Flink aims to support " + "all data types" -> single run
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$body4 = $s4.Shapes.Item(2)
$body4Range = $body4.TextFrame.TextRange
$para1 = $body4Range.Paragraphs(1, 1)
$para1.Text = "TEMP_PLACEHOLDER"
$body4.TextFrame.TextRange.Paragraphs(1, 1).Text = "Flink aims to support all data types"

Write-Host "edits applied"
